$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    $range.Value = "'" + $text
    $range.Style = "Normal"
}

# B2: count 484 -> 482 (numeric)
$ws.Range("B2").Value = 482

# C2: date_min 15440000 -> 1544 (text)
Set-TextValue $ws.Range("C2") "1544"

# B3: count 60 -> 62 (numeric)
$ws.Range("B3").Value = 62

# C4: date_min 15460000 -> 1546 (text)
Set-TextValue $ws.Range("C4") "1546"

# C7/D7: date_min/date_max -> 1548 / 1736 (text)
Set-TextValue $ws.Range("C7") "1548"
Set-TextValue $ws.Range("D7") "1736"

# C11: date_min 16280000 -> 1628 (text)
Set-TextValue $ws.Range("C11") "1628"

# D27: date_max 16780000 -> 1678 (text)
Set-TextValue $ws.Range("D27") "1678"

# C30: date_min 16270000 -> 1627 (text)
Set-TextValue $ws.Range("C30") "1627"

# C40: date_min 15820000 -> 1582 (text)
Set-TextValue $ws.Range("C40") "1582"

# C44/D44: date_min/date_max -> 1641 / 1641 (text)
Set-TextValue $ws.Range("C44") "1641"
Set-TextValue $ws.Range("D44") "1641"

# C57/D57: date_min/date_max -> 00 / 00 (text)
Set-TextValue $ws.Range("C57") "00"
Set-TextValue $ws.Range("D57") "00"

# C61/D61: date_min/date_max -> 1610 / 1610 (text)
Set-TextValue $ws.Range("C61") "1610"
Set-TextValue $ws.Range("D61") "1610"

# C68/D68: date_min/date_max -> 1553 / 1553 (text)
Set-TextValue $ws.Range("C68") "1553"
Set-TextValue $ws.Range("D68") "1553"
